# Update cryptos list values per diff (Mon Mar 11 23:28:01 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.261.70"
$ws.Range("E2").Value = "  +5.22%  "
$ws.Range("D3").Value = "4.075.77"
$ws.Range("E3").Value = "  +5.80%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'521.90"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "'148.69"
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("E7").Value = "  +19.48%  "
$ws.Range("D8").Value = "4.066.87"
$ws.Range("E8").Value = "  +5.77%  "
$ws.Range("D10").Value = "'0.780"
$ws.Range("E10").Value = "  +9.65%  "
$ws.Range("E11").Value = "  +6.45%  "
$ws.Range("D12").Value = "'0.0000334"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "'48.80"
$ws.Range("E13").Value = "  +17.13%  "
$ws.Range("E14").Value = "  +8.85%  "
$ws.Range("D15").Value = "4.719.56"
$ws.Range("E15").Value = "  +5.88%  "
$ws.Range("D16").Value = "4.095.33"
$ws.Range("E16").Value = "  +6.86%  "
$ws.Range("D17").Value = "'14.59"
$ws.Range("E17").Value = "  +4.90%  "
$ws.Range("D18").Value = "'21.42"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'1.24"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D21").Value = "72.363.83"
$ws.Range("E21").Value = "  +5.39%  "
$ws.Range("D22").Value = "'448.55"
$ws.Range("E22").Value = "  +7.07%  "
$ws.Range("D23").Value = "'103.86"
$ws.Range("E23").Value = "  +19.55%  "
$ws.Range("E24").Value = "  +6.44%  "
$ws.Range("D25").Value = "'15.09"
$ws.Range("E25").Value = "  +7.67%  "
$ws.Range("D26").Value = "'4.06"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").Value = "'11.49"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").Value = "'11.14"
$ws.Range("E28").Value = "  +5.55%  "
$ws.Range("D29").Value = "'38.18"
$ws.Range("E29").Value = "  +6.23%  "
$ws.Range("D30").Value = "'5.83"
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("E31").Value = "  +16.59%  "
$ws.Range("D32").Value = "'13.74"
$ws.Range("E32").Value = "  +5.22%  "
$ws.Range("E33").Value = "  +5.17%  "
$ws.Range("D34").Value = "'683.45"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.69"
$ws.Range("E35").Value = "  +13.96%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'67.78"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "'42.50"
$ws.Range("E37").Value = "  +6.81%  "
$ws.Range("D38").Value = "0.0₃0884"
$ws.Range("E38").Value = "  +4.40%  "
$ws.Range("D39").Value = "'0.434"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("E41").Value = "  +10.13%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +5.83%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("E46").Value = "  +13.76%  "
$ws.Range("D47").Value = "'9.92"
$ws.Range("E47").Value = "  +17.46%  "
$ws.Range("D48").Value = "'2.70"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("E50").Value = "  +4.85%  "
$ws.Range("D51").Value = "'0.000284"
$ws.Range("E51").Value = "  +5.12%  "
